$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: locate a province/city row by its name in column A (rows 4-64)
# and overwrite its "Casos totales / Casos activos / Recuperados / Muertes"
# figures with the updated values from the newer data pull.
function Set-CityRow {
    param(
        [string]$Name,
        [int]$CasosTotales,
        [int]$CasosActivos,
        [int]$Recuperados,
        [int]$Muertes
    )
    $cell = $ws.Range("A4:A64").Find($Name, $null, $null, [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
    $r = $cell.Row
    $ws.Cells.Item($r, 2).Value = $CasosTotales
    $ws.Cells.Item($r, 3).Value = $CasosActivos
    $ws.Cells.Item($r, 4).Value = $Recuperados
    $ws.Cells.Item($r, 5).Value = $Muertes
}

# Update timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 12:16"

# Apply the updated case counts for the provinces whose figures changed
Set-CityRow "Madrid"              14597 3031 9741 1825
Set-CityRow "Castilla-La Mancha"  2780  71   2446 263
Set-CityRow "Valencia/Valencia"   1497  23   1413 61
Set-CityRow "Navarra"             1197  23   1141 33
Set-CityRow "La Rioja"            928   43   848  37
Set-CityRow "Alacant/Alicante"    857   12   774  71
Set-CityRow "Cantabria"           510   12   484  14
Set-CityRow "Castello/Castellon"  269   1    257  11
Set-CityRow "Melilla"             38    0    38   0

# The table is kept sorted by "Casos totales" (column B) descending, so
# re-sort the data rows now that the figures above have been refreshed.
$dataRange = $ws.Range("A4:E64")
$sortKey = $ws.Range("B4:B64")
$dataRange.Sort($sortKey, 2, $null, $null, 1, $null, 1, 1)
